$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new negative keywords
$ws.Range("A3").Value = "shit"
$ws.Range("A4").Value = "fuck"

# Bump the sheet's row-outline watermark to 3 without leaving a visible
# grouped row behind: group a throwaway row then delete it so only the
# <sheetFormatPr outlineLevelRow="3"/> high-water mark remains.
$ws.Range("A10").EntireRow.OutlineLevel = 3
$ws.Range("A10").EntireRow.Delete()

# Update selection to match the new target range
$ws.Range("E5:F5").Select()

# Update window dimensions
$excel.ActiveWindow.Width = 28785
$excel.ActiveWindow.Height = 10620
